$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-28 05:08:52"

$wsZhCn.Range("H2").Value = "2016-08-28 05:08:47"
$wsZhCn.Range("K2").Value = "2016-08-28 05:09:08"

$wsDeDe.Range("K2").Value = "2016-08-28 05:09:15"
